# Update "paises.xlsx" COVID-19 country/provincia data (Spain) - 11 Abr 2020, 20:52 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 20:52"

# 2. Update the per-country rows whose figures changed in this refresh.
#    Columns: B=Casos totales C=Nuevos casos D=Casos activos E=Recuperados
#             F=Casos criticos G=Muertes hoy H=Muertes
$updates = @{
    "Estados Unidos"               = @(522320, 19444, 28587, 473647, 10966, 1339, 20086)
    "Suiza"                        = @(25107, 556, 12100, 11971, 386, 34, 1036)
    "Canada"                       = @(23195, 1047, 6309, 16238, 557, 79, 648)
    "Brasil"                       = @(20247, 458, 173, 18984, 296, 22, 1090)
    "Israel"                       = @(10743, 335, 1341, 9301, 175, 6, 101)
    "Irlanda"                      = @(8928, 839, 25, 8583, 194, 33, 320)
    "Noruega"                      = @(6403, 89, 32, 6253, 67, 5, 118)
    "Pakistan"                     = @(5011, 316, 762, 4164, 50, 19, 85)
    "Arabia Saudita"               = @(4033, 382, 720, 3261, 67, 5, 52)
    "Islandia"                     = @(1689, 14, 841, 840, 11, 1, 8)
    "San Martin (Parte Holandesa)" = @(50, 0, 5, 36, 2, 1, 9)
    "Gambia"                       = @(9, 5, 2, 6, 0, 0, 1)
}

$countryRange = $ws.Range("A4:A216")
$cols = @("B", "C", "D", "E", "F", "G", "H")
foreach ($country in $updates.Keys) {
    $cellFound = $countryRange.Find($country)
    $r = $cellFound.Row
    $vals = $updates[$country]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# 3. Countries are kept sorted by "Casos totales" (column B) descending, so
#    re-sort the whole data block now that the figures above moved some
#    countries past their neighbours.
$dataRange = $ws.Range("A4:H216")
$dataRange.Sort($ws.Range("B4:B216"), 2, $null, $null, 1, $null, 1, 2)

# 4. A handful of countries tie exactly on every figure with their neighbour
#    after the refresh; the source feed listed them in swapped order, so
#    mirror that for the tied pairs (values are identical, only the two
#    labels trade places).
$tiedSwapPairs = @(
    @("Sudan", "Angola"),
    @("Republica de Africa Central", "Sierra Leona"),
    @("Nicaragua", "Islas Turcas y Caicos"),
    @("Sahara Occidental", "Santo Tome y Principe"),
    @("Papua Nueva Guinea", "Bonaire, San Eustaquio y Saba"),
    @("San Pedro y Miquelon", "Yemen")
)
$fullRange = $ws.Range("A4:A216")
$allCols = @("A", "B", "C", "D", "E", "F", "G", "H")
foreach ($pair in $tiedSwapPairs) {
    $r1 = $fullRange.Find($pair[0]).Row
    $r2 = $fullRange.Find($pair[1]).Row
    foreach ($col in $allCols) {
        $v1 = $ws.Range($col + $r1).Value2
        $v2 = $ws.Range($col + $r2).Value2
        $ws.Range($col + $r1).Value = $v2
        $ws.Range($col + $r2).Value = $v1
    }
}
